# Updated Simulated Annealing route optimization results (Column C: Customer)
# Fixes an issue where the SA algorithm wasn't accepting worse moves,
# so the "Customer" assignment per sequence position has changed for
# every non-starting row in each paperboy's route.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$customerUpdates = @{
    3 = 55
    4 = 69
    5 = 73
    6 = 74
    7 = 103
    8 = 110
    9 = 115
    10 = 120
    11 = 109
    12 = 102
    13 = 112
    14 = 97
    15 = 76
    16 = 66
    17 = 57
    18 = 58
    19 = 48
    20 = 51
    21 = 2
    22 = 5
    24 = 33
    25 = 28
    26 = 23
    27 = 16
    28 = 12
    29 = 11
    30 = 8
    31 = 13
    32 = 30
    34 = 52
    35 = 50
    36 = 45
    37 = 34
    38 = 35
    39 = 29
    40 = 20
    41 = 19
    42 = 31
    43 = 25
    44 = 15
    45 = 4
    46 = 6
    47 = 14
    48 = 38
    49 = 59
    50 = 62
    51 = 67
    52 = 68
    53 = 81
    54 = 101
    55 = 98
    56 = 85
    57 = 114
    58 = 104
    59 = 96
    60 = 86
    61 = 82
    62 = 77
    63 = 71
    65 = 83
    66 = 44
    67 = 47
    68 = 39
    69 = 36
    70 = 32
    71 = 24
    72 = 22
    73 = 10
    74 = 1
    75 = 7
    76 = 46
    77 = 43
    78 = 84
    79 = 99
    80 = 108
    81 = 117
    82 = 91
    83 = 105
    84 = 106
    86 = 54
    87 = 87
    88 = 93
    89 = 94
    90 = 80
    91 = 78
    92 = 70
    93 = 64
    94 = 37
    96 = 89
    97 = 100
    98 = 107
    99 = 111
    100 = 116
    101 = 119
    102 = 118
    103 = 113
    104 = 61
    105 = 60
    106 = 65
    107 = 88
    108 = 95
    109 = 92
    110 = 90
    111 = 63
    112 = 79
    113 = 56
    114 = 49
    115 = 42
    116 = 18
    117 = 21
    118 = 40
    119 = 75
    120 = 53
    121 = 41
    122 = 27
    123 = 17
    124 = 3
    125 = 9
}

foreach ($row in $customerUpdates.Keys) {
    $ws.Cells.Item($row, 3).Value = $customerUpdates[$row]
}

Write-Host "Updated $($customerUpdates.Count) Customer values in column C."
